$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "No. Seri"
$ws.Range("A1").Select()
